$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (A1/B1) already exists ("単語" / "意味").
# Add the new word entry in row 2, column A ("さいけいせい").
$ws.Range("A2").Value = "さいけいせい"

# Excel auto-sized column A to fit the new (wider) Japanese content.
$ws.Columns("A").ColumnWidth = 11.8

# Reflect the cursor resting on the next empty row after data entry.
$ws.Range("A3").Select() | Out-Null
